# Insert a new row at row 398 (this shifts rows 398:509 down to 399:510,
# extending the used range to A1:T510) and populate it with the new data
# record, while leaving the columns that remain identical to the old row 398
# (A,B,C,E,F,G,H,I,J,K,L,Q,T) copied over automatically by the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A398").EntireRow.Insert()

$ws.Range("A398").Value = 11
$ws.Range("B398").Value = 'Vega Monumental Concepción'
$ws.Range("C398").Value = 'Bíobío'
$ws.Range("D398").Value = 44736
$ws.Range("E398").Value = 8
$ws.Range("F398").Value = 'Fruta'
$ws.Range("G398").Value = 100102
$ws.Range("H398").Value = 'Cítricos'
$ws.Range("I398").Value = 100102003
$ws.Range("J398").Value = 'Limón'
$ws.Range("K398").Value = 'Sin especificar'
$ws.Range("L398").Value = '1a amarillo'
$ws.Range("M398").Value = 350
$ws.Range("N398").Value = 4500
$ws.Range("O398").Value = 5000
$ws.Range("P398").Value = 4714
$ws.Range("Q398").Value = '$/malla 16 kilos'
$ws.Range("R398").Value = 'Provincia de Quillota'
$ws.Range("S398").Value = 295
$ws.Range("T398").Value = 16
